$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 14 (pushes the existing blank row13 down,
# along with the totals row and the two merged "spacer" rows below it).
# After two inserts at the same spot:
#   row13 -> still the original blank row
#   row14 -> new blank row (copy of row13 formatting)
#   row15 -> original blank row13 (shifted down)
#   row16 -> totals row (was row14)
#   row17 -> spacer row (was row15)
#   row18 -> spacer row (was row16)
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

# Fill in the "Day 10" entry into row 13
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "24/5/2024"
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = "Finished Update product + delete product + delete product by admin"

# Fill in the "Day 11" entry into row 14
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "25/5/2024"
$ws.Range("C14").Value = 2.5
$ws.Range("D14").Value = "Finished User Ban / Remove Ban functionality"

# Update the total formula to include the two new rows of hours (rows 4-14)
$ws.Range("D16").Formula = "=SUM(C4:C14)"

# Setting the formula can make Excel mark the row with an explicit custom
# height; auto-fit it back so the row keeps its default height.
$ws.Rows.Item(16).AutoFit()

# Match the saved selection state from the edited workbook
$ws.Range("D16:D18").Select()
